# Correcion error al hacer clic en datagridview
#
# - Remove the two "Da error al hacer clic en el datagridview" task-notes
#   (column B) on the "Compra nueva" / "Venta nueva" rows, and flip those
#   two rows (and the now-resolved "Poder editar la fecha de una Compra"
#   row) from pending/flagged to done (green).
# - Remove the now-empty spacer column C (column D shifts left into C).
# - Re-order a few finished/pending items in the (new) column C and add a
#   brand new pending item: "Cambio de Precios Costos y Utilidades".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Drop the resolved "Da error al hacer clic en el datagridview" notes
#    and mark their rows (and the related date-editing fix) as done.
# ---------------------------------------------------------------------
$ws.Range("B19").Clear()
$ws.Range("B24").Clear()

$ws.Range("A20").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Remove the blank spacer column C -- column D (and E) shift left.
# ---------------------------------------------------------------------
$ws.Columns("C:C").Delete()

# ---------------------------------------------------------------------
# 3) Re-order / insert items in the (now) column C, rows 20-32.
#    Copy formatting first (values don't matter for a format copy), then
#    overwrite the text.
# ---------------------------------------------------------------------
$ws.Range("C25").Copy()
$ws.Range("C31:C32").PasteSpecial(-4122)

$ws.Range("C19").Copy()
$ws.Range("C20:C22").PasteSpecial(-4122)
$ws.Range("C24").PasteSpecial(-4122)

$ws.Range("A29").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C20").Value = "Acerca de"
$ws.Range("C21").Value = "Ayuda"
$ws.Range("C22").Value = "En todos los ""Ver un"" poder cambiar el ID para buscar otro"
$ws.Range("C23").Value = "Cambio de Precios Costos y Utilidades"
$ws.Range("C24").Value = "Asegurarse que no da error al querer borrar un registro que tenga relaciones"
$ws.Range("C25").Value = "Listado de Clientes"
$ws.Range("C26").Value = "Listado de Productos"
$ws.Range("C27").Value = "Listado de Proveedores"
$ws.Range("C28").Value = "Reporte de Inventario"
$ws.Range("C29").Value = "Reporte de Utilidad"
$ws.Range("C30").Value = "Reporte de Ventas"
$ws.Range("C31").Value = "Seguridad de la Aplicaciòn"
$ws.Range("C32").Value = "Instalador"

# ---------------------------------------------------------------------
# 4) Restore the selection/view to match the edited area.
# ---------------------------------------------------------------------
$ws.Range("B24").Select()
